$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 34: cohort 2020, period_index 2 -> num_customers 80 -> 81
$ws.Range("C34").Value2 = 81
$c34 = $ws.Range("C34").Value2
$d34 = $ws.Range("D34").Value2
$ws.Range("E34").Value2 = $c34 / $d34

# Row 37: cohort 2023, period_index 0 -> num_customers/cohort_size 883 -> 889
$ws.Range("C37").Value2 = 889
$ws.Range("D37").Value2 = 889
